# Restructure the "Input" sheet to the standard template column layout
# (A:발주일자, B:납기일자, C:거래처명, D:거래처 이메일, E:납품처명, F:납품처 이메일,
#  G:프로젝트명, H:대분류, I:중분류, J:소분류, K:품목명, L:규격, M:수량, N:단가,
#  O:총금액, P:비고) and drop the old 17th column (old Q / 비고), and also drop the
# leftover empty "비고" cells in the 갑지/을지 summary sheets' I column.

$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("Input")

# Wipe all existing content/formatting (old A:Q layout, including the bold/
# bordered header style) so the sheet can be rebuilt from scratch.
$wsInput.Cells.Clear()

# --- Header row -------------------------------------------------------
$headers = @("발주일자","납기일자","거래처명","거래처 이메일","납품처명","납품처 이메일","프로젝트명","대분류","중분류","소분류","품목명","규격","수량","단가","총금액","비고")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $wsInput.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- Data rows ---------------------------------------------------------
# Columns A, B hold plain ISO date-strings that must stay text (Excel would
# otherwise silently coerce "2025-08-20" into a date serial number), so the
# whole A:B block is forced to Text format before the values are poured in,
# then restored to the default "Normal" style so no stray numFmt survives.
$dataRows = @(
    @("2025-08-20","2025-09-04","제이비엔지니어링","제이비엔지니어링@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","1. 원자재비","4) ALUM. 창호","A. 압출","5월 청구분","KS규격-1",46,4910,248446),
    @("2025-09-04","2025-09-11","제이비엔지니어링","제이비엔지니어링@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","5. 운반비","일반자재","기타","5월 운반비","KS규격-2",1,0,0),
    @("2025-08-31","2025-10-06","제이비엔지니어링","제이비엔지니어링@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","1. 원자재비","4) ALUM. 창호","A. 압출","IJ-15861","KS규격-3",1,458040,503844),
    @("2025-09-08","2025-09-04","제이비엔지니어링","제이비엔지니어링@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","1. 원자재비","4) ALUM. 창호","A. 압출","4월 청구분","KS규격-4",519,4950,2825955)
)

$wsInput.Range("A2:B5").NumberFormat = "@"

for ($r = 0; $r -lt $dataRows.Count; $r++) {
    $rowValues = $dataRows[$r]
    $rowNum = $r + 2
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $wsInput.Cells.Item($rowNum, $c + 1).Value = $rowValues[$c]
    }
    # Column P (비고) is intentionally left untouched/empty for every data row.
}

$wsInput.Range("A2:B5").Style = "Normal"

# --- 갑지 / 을지 summary sheets ----------------------------------------
# The trailing empty "비고" (remark) cells in column I, rows 2-5, are removed
# entirely rather than kept as empty strings.
foreach ($sheetName in @("갑지", "을지")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("I2:I5").ClearContents()
}
